# Auto-generated edit script: applies numeric value updates described in the commit diff
# across sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR (columns H-N, per-row P&L figures).
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 5143
$ws.Range("I18").Value = 1095
$ws.Range("J18").Value = 6077.154
$ws.Range("K18").Value = 1095
$ws.Range("L18").Value = 6077.154
$ws.Range("M18").Value = -811
$ws.Range("N18").Value = -6645.154
$ws.Range("H20").Value = 4372.25
$ws.Range("I20").Value = 4372.25
$ws.Range("K20").Value = 4372.25
$ws.Range("M20").Value = -4142.25
$ws.Range("H35").Value = 4372.25
$ws.Range("I35").Value = 4372.25
$ws.Range("K35").Value = 4372.25
$ws.Range("M35").Value = -3993.25
$ws.Range("H80").Value = 13894126
$ws.Range("I80").Value = 25001228
$ws.Range("J80").Value = 10249.75
$ws.Range("K80").Value = 75003684
$ws.Range("L80").Value = 30749.25
$ws.Range("M80").Value = -75002686
$ws.Range("N80").Value = -32745.25
$ws.Range("H83").Value = 13894126
$ws.Range("I83").Value = 25001228
$ws.Range("J83").Value = 10249.75
$ws.Range("K83").Value = 225011052
$ws.Range("L83").Value = 92247.75
$ws.Range("M83").Value = -225006060
$ws.Range("N83").Value = -102231.75
$ws.Range("H98").Value = 1194.4706
$ws.Range("I98").Value = 1194.4706
$ws.Range("K98").Value = 1194.4706
$ws.Range("M98").Value = 303.5293999999999
$ws.Range("H111").Value = 4601.8335
$ws.Range("I111").Value = 2632.7144
$ws.Range("J111").Value = 7358.6
$ws.Range("K111").Value = 7898.1432
$ws.Range("L111").Value = 22075.8
$ws.Range("M111").Value = -4831.1432
$ws.Range("N111").Value = -28209.8
$ws.Range("H121").Value = 11020.6
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 11020.6
$ws.Range("K121").Value = 0
$ws.Range("L121").Value = 33061.8
$ws.Range("M121").ClearContents()
$ws.Range("N121").Value = -36555.8
$ws.Range("H122").Value = 1194.4706
$ws.Range("I122").Value = 1194.4706
$ws.Range("K122").Value = 3583.4118
$ws.Range("M122").Value = -1133.4118
$ws.Range("H132").Value = 2723.1428
$ws.Range("I132").Value = 1394.6154
$ws.Range("K132").Value = 4183.8462
$ws.Range("M132").Value = -1653.8462
$ws.Range("H138").Value = 7099.1304
$ws.Range("I138").Value = 11065.782
$ws.Range("K138").Value = 33197.346
$ws.Range("M138").Value = -28057.346

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8334167
$ws.Range("I32").Value = 8772692
$ws.Range("K32").Value = 8772692
$ws.Range("M32").Value = -8772405
$ws.Range("H74").Value = 5037.6523
$ws.Range("I74").Value = 1921.6666
$ws.Range("J74").Value = 10880.125
$ws.Range("K74").Value = 1921.6666
$ws.Range("L74").Value = 10880.125
$ws.Range("M74").Value = -1047.6666
$ws.Range("N74").Value = -12628.125
$ws.Range("H77").Value = 5037.6523
$ws.Range("I77").Value = 1921.6666
$ws.Range("J77").Value = 10880.125
$ws.Range("K77").Value = 9608.333000000001
$ws.Range("L77").Value = 54400.625
$ws.Range("M77").Value = -5240.333000000001
$ws.Range("N77").Value = -63136.625
$ws.Range("H122").Value = 3585.1667
$ws.Range("I122").Value = 2286.8
$ws.Range("K122").Value = 6860.400000000001
$ws.Range("M122").Value = -4410.400000000001

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 674.75
$ws.Range("I22").Value = 674.75
$ws.Range("K22").Value = 674.75
$ws.Range("M22").Value = -501.75
$ws.Range("H60").Value = 0
$ws.Range("J60").Value = 0
$ws.Range("L60").Value = 0
$ws.Range("N60").ClearContents()
$ws.Range("H134").Value = 2925.0142
$ws.Range("I134").Value = 1170.0741
$ws.Range("J134").Value = 8499.529
$ws.Range("K134").Value = 3510.2223
$ws.Range("L134").Value = 25498.587
$ws.Range("M134").Value = -975.2223000000004
$ws.Range("N134").Value = -30568.587
$ws.Range("H141").Value = 70000
$ws.Range("J141").Value = 70000
$ws.Range("L141").Value = 70000
$ws.Range("N141").Value = -80360

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 83586.64999999999
$ws.Range("J31").Value = 22852.684
$ws.Range("L31").Value = 22852.684
$ws.Range("N31").Value = -23442.684
$ws.Range("H34").Value = 83586.64999999999
$ws.Range("J34").Value = 22852.684
$ws.Range("L34").Value = 22852.684
$ws.Range("N34").Value = -23256.684
$ws.Range("H58").Value = 3868.3333
$ws.Range("I58").Value = 2061.625
$ws.Range("J58").Value = 5313.7
$ws.Range("K58").Value = 2061.625
$ws.Range("L58").Value = 5313.7
$ws.Range("M58").Value = -1858.625
$ws.Range("N58").Value = -5719.7
$ws.Range("H62").Value = 5562.875
$ws.Range("I62").Value = 2417.8333
$ws.Range("K62").Value = 2417.8333
$ws.Range("M62").Value = -1793.8333
$ws.Range("H65").Value = 5562.875
$ws.Range("I65").Value = 2417.8333
$ws.Range("K65").Value = 12089.1665
$ws.Range("M65").Value = -8969.166499999999
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()
$ws.Range("H136").Value = 3868.3333
$ws.Range("I136").Value = 2061.625
$ws.Range("J136").Value = 5313.7
$ws.Range("K136").Value = 6184.875
$ws.Range("L136").Value = 15941.1
$ws.Range("M136").Value = -3634.875
$ws.Range("N136").Value = -21041.1

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 1330.55
$ws.Range("I129").Value = 649.7143
$ws.Range("J129").Value = 2919.1667
$ws.Range("K129").Value = 1949.1429
$ws.Range("L129").Value = 8757.500100000001
$ws.Range("M129").Value = 3050.8571
$ws.Range("N129").Value = -18757.5001
$ws.Range("H140").Value = 2139.8215
$ws.Range("I140").Value = 1735.4348
$ws.Range("K140").Value = 5206.3044
$ws.Range("M140").Value = -26.30439999999999
$ws.Range("H141").Value = 2546
$ws.Range("I141").Value = 2546
$ws.Range("K141").Value = 7638
$ws.Range("M141").Value = -2458

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H23").Value = 2514.2856
$ws.Range("I23").Value = 300
$ws.Range("J23").Value = 3400
$ws.Range("K23").Value = 300
$ws.Range("L23").Value = 3400
$ws.Range("M23").Value = -77
$ws.Range("N23").Value = -3846
$ws.Range("H113").Value = 1940.4166
$ws.Range("I113").Value = 1737.7
$ws.Range("K113").Value = 1737.7
$ws.Range("M113").Value = 432.3

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 55862.168
$ws.Range("I136").Value = 1757.5454
$ws.Range("K136").Value = 5272.6362
$ws.Range("M136").Value = -2722.6362

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1897.9615
$ws.Range("I107").Value = 997.6667
$ws.Range("J107").Value = 3923.625
$ws.Range("K107").Value = 2993.0001
$ws.Range("L107").Value = 11770.875
$ws.Range("M107").Value = -1073.0001
$ws.Range("N107").Value = -15610.875
$ws.Range("H116").Value = 129777.375
$ws.Range("J116").Value = 129777.375
$ws.Range("L116").Value = 129777.375
$ws.Range("N116").Value = -138955.375

